$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/medicare-coverage"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": clear the duplicated constraint text that used to
#     sit on the root "Extension" row (row 2, column AI = "Constraint(s)").
#     That same text already lives on the "Extension.extension" row (row 4).
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
